$wb = $excel.ActiveWorkbook

# --- Rename the first worksheet (drop the mangled suffix) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "0105"

# Fix up the Print_Titles defined name so it quotes the (now short) sheet name
foreach ($n in $wb.Names) {
    if ($n.Name -like "*Print_Titles*") {
        $n.RefersTo = "='0105'!`$1:`$1"
    }
}

# --- Add a new worksheet "Sheet1" right after "0105" ---
$newWs = $wb.Worksheets.Add($null, $ws1)
$newWs.Name = "Sheet1"

# Bring over the header row + first data row formatting from "0105"
$ws1.Range("A1:T1").Copy()
$newWs.Range("A1:T1").PasteSpecial(-4122)
$ws1.Range("A2:T2").Copy()
$newWs.Range("A2:T2").PasteSpecial(-4122)

$newWs.Rows.Item(1).RowHeight = 12
$newWs.Rows.Item(2).RowHeight = 49.5

# Header row (same labels as sheet "0105")
$newWs.Range("A1").Value = "INV:"
$newWs.Range("B1").Value = "ItemCode"
$newWs.Range("C1").Value = "CO"
$newWs.Range("D1").Value = "CartonNo"
$newWs.Range("E1").Value = "PalletNo"
$newWs.Range("F1").Value = "PartNumber"
$newWs.Range("G1").Value = "DeliveryQty"
$newWs.Range("H1").Value = "Mfgr"
$newWs.Range("I1").Value = "PO[1]"
$newWs.Range("J1").Value = "POQty[1]"
$newWs.Range("K1").Value = "PO[2]"
$newWs.Range("L1").Value = "POQty[2]"
$newWs.Range("M1").Value = "PO[3]"
$newWs.Range("N1").Value = "POQty[3]"
$newWs.Range("O1").Value = "PO[4]"
$newWs.Range("P1").Value = "POQty[4]"
$newWs.Range("Q1").Value = "PO[5]"
$newWs.Range("R1").Value = "POQty[5]"

# Data row
$newWs.Range("A2").Value = "K-1-5-15A"
$newWs.Range("B2").Value = "7YCS12B1008+H01"
$newWs.Range("C2").Value = "INDONESIA"
$newWs.Range("D2").Value = "50-51"
$newWs.Range("E2").Value = 2
$newWs.Range("F2").Value = "50501201R1F"
$newWs.Range("G2").Value = 2800
$newWs.Range("H2").Value = "JSTT00"
$newWs.Range("I2").Value = 10746487
$newWs.Range("J2").Value = 80
$newWs.Range("K2").Value = 10748713
$newWs.Range("L2").Value = 2720

# Selections: "Sheet1" remembers D9, then focus returns to "0105" at D8
$newWs.Range("D9").Select() | Out-Null
$ws1.Select() | Out-Null
$ws1.Range("D8").Select() | Out-Null
